# Generate Report for handoff
# Update the "b.md.md" row (row 3) across the Overview, zh-cn and de-de
# sheets to reflect that a new handoff has been generated.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# ---- zh-cn sheet ---------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-02-15 02:59:34"

foreach ($h in $zhcn.Hyperlinks) {
    if ($h.Range.Address() -eq '$C$3') {
        $h.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
    }
}

# ---- de-de sheet ---------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$dede.Range("D3").Value = "2016-02-15 02:59:47"

foreach ($h in $dede.Hyperlinks) {
    if ($h.Range.Address() -eq '$C$3') {
        $h.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
    }
}
